# Update the two data rows (row 2 and row 3) with their new values,
# matching the reformatted figures produced after the "black" reformat.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -0.01132772767135682
$ws.Range("D2").Value = 0.1329501444668215
$ws.Range("E2").Value = 0.2452044431023112
$ws.Range("F2").Value = 0.3210583506375297
$ws.Range("G2").Value = 0.3622452035370414
$ws.Range("H2").Value = 0.1654327651585601
$ws.Range("I2").Value = 0.2818359714947122
$ws.Range("J2").Value = 0.355822772016644
$ws.Range("K2").Value = 0.4034558434828278

$ws.Range("C3").Value = 0.04167199304037311
$ws.Range("D3").Value = 0.1945323064040473
$ws.Range("E3").Value = 0.315326141460165
$ws.Range("F3").Value = 0.3952634884591798
$ws.Range("G3").Value = 0.4363374629815695
$ws.Range("H3").Value = 0.2269676979385108
$ws.Range("I3").Value = 0.3563520789343781
$ws.Range("J3").Value = 0.4342047310983628
$ws.Range("K3").Value = 0.4774351088200124
